$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column M (rows 4-17) into the new column N, one row
# at a time, so each new cell inherits the same style as its neighbour in M.
for ($r = 4; $r -le 17; $r++) {
    $ws.Range("M$r").Copy()
    $ws.Range("N$r").PasteSpecial(-4122)  # xlPasteFormats
}

# New "2020" data column values (row 15 is a blank spacer row with no value).
$ws.Range("N4").Value  = 2020
$ws.Range("N5").Value  = 11.4
$ws.Range("N6").Value  = 14.7
$ws.Range("N7").Value  = 9
$ws.Range("N8").Value  = 10.8
$ws.Range("N9").Value  = 4.7
$ws.Range("N10").Value = 5.0999999999999996
$ws.Range("N11").Value = 3.4
$ws.Range("N12").Value = 19.7
$ws.Range("N13").Value = 18.8
$ws.Range("N14").Value = 6.8
$ws.Range("N16").Value = 12.5
$ws.Range("N17").Value = 10.7

# Update the view: scroll so column E is the leftmost visible column, and
# move the selection to S18.
$ws.Range("S18").Select()
